$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title 1 shape: "A" + " " + "slide" -> single run "A slide"
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "."
$titleRange.Text = "A slide"

# Table cell (row 1, col 2): "a" + " " + "table" -> single run "a table"
$tableShape = $s.Shapes.Item(3)
$cellRange = $tableShape.Table.Cell(1, 2).Shape.TextFrame.TextRange
$cellRange.Text = "."
$cellRange.Text = "a table"

# TextBox 3: "Plus" + " " + "an" + " " + "image" -> single run "Plus an image"
$textBoxRange = $s.Shapes.Item(7).TextFrame.TextRange
$textBoxRange.Text = "."
$textBoxRange.Text = "Plus an image"
